$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '54.615.75'
$ws.Range("E2").Value = '  -6.44%  '

$ws.Range("D3").Value = '2.443.62'
$ws.Range("E3").Value = '  -9.35%  '

$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '467.97'
$ws.Range("E5").Value = '  -6.38%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '131.32'
$ws.Range("E6").Value = '  -5.54%  '

$ws.Range("E7").Value = '  +0.32%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.494'
$ws.Range("E8").Value = '  -5.89%  '

$ws.Range("D9").Value = '2.448.21'
$ws.Range("E9").Value = '  -9.57%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0952'
$ws.Range("E10").Value = '  -8.74%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.31'
$ws.Range("E11").Value = '  -11.95%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.316'
$ws.Range("E12").Value = '  -8.82%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.121'
$ws.Range("E13").Value = '  -3.92%  '

$ws.Range("D14").Value = '2.877.64'
$ws.Range("E14").Value = '  -9.17%  '

$ws.Range("D15").Value = '54.722.61'
$ws.Range("E15").Value = '  -6.35%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000134'
$ws.Range("E16").Value = '  +0.30%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '19.67'
$ws.Range("E17").Value = '  -7.96%  '

$ws.Range("D18").Value = '2.450.01'
$ws.Range("E18").Value = '  -9.40%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.23'
$ws.Range("E19").Value = '  -10.21%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '312.02'
$ws.Range("E20").Value = '  -6.37%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.55'
$ws.Range("E21").Value = '  -12.48%  '

$ws.Range("E22").Value = '  +0.50%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.67'
$ws.Range("E23").Value = '  +0.75%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.37'
$ws.Range("E24").Value = '  -13.32%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '56.43'
$ws.Range("E25").Value = '  -10.09%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.01'
$ws.Range("E26").Value = '  +1.19%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.157'
$ws.Range("E27").Value = '  -7.85%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.386'
$ws.Range("E28").Value = '  -8.66%  '

$ws.Range("D29").Value = '2.558.77'
$ws.Range("E29").Value = '  -9.29%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.13'
$ws.Range("E30").Value = '  -3.53%  '

$ws.Range("E31").Value = '  +0.05%  '

$ws.Range("D32").Value = '0.0₃0715'
$ws.Range("E32").Value = '  -12.31%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '145.96'
$ws.Range("E33").Value = '  -3.05%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '17.78'
$ws.Range("E34").Value = '  -6.67%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.43'
$ws.Range("E35").Value = '  -9.65%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.02'
$ws.Range("E36").Value = '  -6.10%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.58'
$ws.Range("E37").Value = '  -13.38%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.05'
$ws.Range("E38").Value = '  -5.37%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.791'
$ws.Range("E39").Value = '  -14.70%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").Value = '  +0.26%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '32.93'
$ws.Range("E41").Value = '  -6.78%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.600'
$ws.Range("E42").Value = '  +1.15%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0524'
$ws.Range("E43").Value = '  -5.49%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.26'
$ws.Range("E44").Value = '  -7.97%  '

$ws.Range("B45").Value = 'WhiteBITCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.11'
$ws.Range("E45").Value = '  -2.45%  '

$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.24'
$ws.Range("E46").Value = '  -9.31%  '

$ws.Range("D47").Value = '1.941.40'
$ws.Range("E47").Value = '  -10.68%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0885'
$ws.Range("E48").Value = '  +0.20%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0217'
$ws.Range("E49").Value = '  -3.45%  '

$ws.Range("B50").Value = 'Bittensor'
$ws.Range("C50").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '233.08'
$ws.Range("E50").Value = '  +6.44%  '

$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.21'
$ws.Range("E51").Value = '  -8.42%  '
